$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rich-text shared-string edits -----------------------------------------
# A8: "Volume 31   Number  2" -> "...3" (edit just the trailing run's text,
# then restore that run's font so the edited digit keeps the original look).
$a8 = $ws.Range("A8")
$a8.Characters(21, 1).Text = "3"
$a8run = $a8.Characters(21, 1)
$a8run.Font.Size = 10
$a8run.Font.Name = "Andale WT"

# C9: "Report Covering the Week  1/8/2024  Through  1/14/2024"
#     -> "...1/15/2024  Through  1/21/2024"
$c9 = $ws.Range("C9")
$c9.Characters(27, 8).Text = "1/15/2024"
$c9run1 = $c9.Characters(27, 9)
$c9run1.Font.Size = 10
$c9run1.Font.Name = "Andale WT"
$c9.Characters(47, 9).Text = "1/21/2024"
$c9run2 = $c9.Characters(47, 9)
$c9run2.Font.Size = 10
$c9run2.Font.Name = "Andale WT"

# --- Plain numeric value updates -------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 54.545454545454
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = 18.181818181818
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -27.777777777777
$ws.Range("N16").Value = -90.441176470588
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -57.142857142857
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 16
$ws.Range("J17").Value = 21
$ws.Range("K17").Value = -23.809523809523
$ws.Range("L17").Value = -27.272727272727
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -60.975609756097
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -64.705882352941
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 12
$ws.Range("K18").Value = -50
$ws.Range("L18").Value = -25
$ws.Range("M18").Value = -50
$ws.Range("N18").Value = -97.183098591549
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -62.264150943396
$ws.Range("I19").Value = 17
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = -57.5
$ws.Range("L19").Value = -52.777777777777
$ws.Range("M19").Value = -55.263157894736
$ws.Range("N19").Value = -62.222222222222
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 8
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 33.333333333333
$ws.Range("M20").Value = -42.857142857142
$ws.Range("N20").Value = -95.294117647058
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -42.307692307692
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = -34.710743801652
$ws.Range("I21").Value = 60
$ws.Range("J21").Value = 92
$ws.Range("K21").Value = -34.782608695652
$ws.Range("L21").Value = -31.034482758620
$ws.Range("M21").Value = -39.393939393939
$ws.Range("N21").Value = -90.259740259740
$ws.Range("F22").Value = 2
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = -20.512820512820
$ws.Range("F24").Value = 127
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = 11.403508771929
$ws.Range("I24").Value = 104
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = 30
$ws.Range("L24").Value = 48.571428571428
$ws.Range("M24").Value = 36.842105263157
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = 9.090909090909
$ws.Range("I25").Value = 32
$ws.Range("J25").Value = 31
$ws.Range("K25").Value = 3.225806451612
$ws.Range("L25").Value = 28
$ws.Range("M25").Value = -15.789473684210
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 2
$ws.Range("K26").Value = -50
$ws.Range("L26").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -50
$ws.Range("H28").Value = -100
$ws.Range("H29").Value = -100

# --- Cells changing from number <-> text (shared-string) type --------------
# Row 22: G22/H22 become text ("0" / "***.*") formatted like the other text
# cells in the row (style of C22).
$ws.Range("G22").Value = "'0"
$ws.Range("H22").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").PasteSpecial(-4122)

# Row 27: C27/I27 become numbers (style of J27, a numeric cell in same row).
$ws.Range("C27").Value = 1
$ws.Range("I27").Value = 1
$ws.Range("J27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("I27").PasteSpecial(-4122)

# Rows 28/29: F28/F29 become text ("0"), formatted like the other text cells
# in the row (style of C28 / C29).
$ws.Range("F28").Value = "'0"
$ws.Range("C28").Copy()
$ws.Range("F28").PasteSpecial(-4122)

$ws.Range("F29").Value = "'0"
$ws.Range("C29").Copy()
$ws.Range("F29").PasteSpecial(-4122)
